# "Now using faster multiplier in cel to vera"
# Insert a new row for "Multiplier Output" in the Vera sheet, right above
# the existing "Vacant" row, and shrink the old Vacant block to make room
# for it (2 + 126 = 128, same total as before).

$wb = $excel.ActiveWorkbook
$vera = $wb.Worksheets.Item("Vera")
$sprites = $wb.Worksheets.Item("Sprite Addresses")

# Insert a new row 4 (pushes old row4 "Vacant" etc. down to row5+)
$vera.Rows.Item(4).Insert()

# Populate the new row 4 with "Multiplier Output" data
$vera.Range("A4").Value = "Multiplier Output"
$vera.Range("B4").Formula = "=B3+ C3"
$vera.Range("C4").Value = 2
$vera.Range("D4").Formula = "=DEC2HEX(B3 + C3)"
$vera.Range("E4").Formula = "=DEC2HEX(HEX2DEC(D4)+C4-1)"

# Fix up the old "Vacant" row (now row 5): shrink size from 128 to 126,
# and correct its formulas to reference the new row 4 above it.
$vera.Range("B5").Formula = "=B4+ C4"
$vera.Range("C5").Value = 126
$vera.Range("D5").Formula = "=DEC2HEX(B4 + C4)"

# Update selection on Vera sheet and make it the active/selected tab
$sprites.Range("D1").Select()
$vera.Activate()
$vera.Range("F8").Select()
